$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.744.50'
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('D3').Value = '2.085.81'
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').Value = "'0.616"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = "'60.20"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('D10').Value = "'0.0843"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '2.395.01'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('E13').Value = '  +4.13%  '
$ws.Range('D14').Value = "'21.86"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = "'0.798"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.68%  '
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').Value = '2.082.65'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').Value = '38.697.07'
$ws.Range('D19').Value = "'71.45"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.99%  '
$ws.Range('D20').Value = "'6.03"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.92%  '
$ws.Range('D21').Value = '0.0₃0840'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('D22').Value = "'227.17"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').Value = "'2.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('D26').Value = "'170.93"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('E28').Value = '  +7.14%  '
$ws.Range('E29').Value = '  +12.65%  '
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  +5.01%  '
$ws.Range('E33').Value = '  +2.75%  '
$ws.Range('E34').Value = '  +3.53%  '
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = "'17.89"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.30%  '
$ws.Range('D41').Value = "'0.0227"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.49%  '
$ws.Range('D42').Value = '1.541.15'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').Value = "'100.49"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.57%  '
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').Value = "'7.70"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.00%  '
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').Value = "'4.12"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '2.283.62'
$ws.Range('E51').Value = '  +2.07%  '
